# Update column C ("Förändrad") for rows 2-18 from serial date 45190
# (2023-09-21) to 45192 (2023-09-23), keeping existing date formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45190) {
        $cell.Value = 45192
    }
}
